$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B7").Value = 6.0
$ws.Range("C7").Value = 3.0350000858306885
$ws.Range("D7").Value = 140.0
